$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1992619926199262
$ws.Range("C2").Value = 0.5645756457564576
$ws.Range("J2").Value = 0.01476014760147601
$ws.Range("P2").Value = 0.1439114391143911
$ws.Range("S2").Value = 0.07749077490774908
$ws.Range("B3").Value = 0.01282051282051282
$ws.Range("J3").Value = 0.02564102564102564
$ws.Range("P3").Value = 0.7756410256410257
$ws.Range("S3").Value = 0.1858974358974359
$ws.Range("J4").Value = 0.07317073170731707
$ws.Range("P4").Value = 0.6097560975609756
$ws.Range("S4").Value = 0.3170731707317073
$ws.Range("B6").Value = 0.03508771929824561
$ws.Range("D6").Value = 0.02192982456140351
$ws.Range("F6").Value = 0.07017543859649122
$ws.Range("J6").Value = 0.3421052631578947
$ws.Range("O6").Value = 0.01754385964912281
$ws.Range("Q6").Value = 0.1491228070175439
$ws.Range("R6").Value = 0.06140350877192982
$ws.Range("S6").Value = 0.3026315789473684
$ws.Range("B7").Value = 0.07253886010362694
$ws.Range("D7").Value = 0.02072538860103627
$ws.Range("E7").Value = 0.005181347150259068
$ws.Range("F7").Value = 0.06217616580310881
$ws.Range("J7").Value = 0.1450777202072539
$ws.Range("O7").Value = 0.02072538860103627
$ws.Range("Q7").Value = 0.1191709844559585
$ws.Range("R7").Value = 0.1295336787564767
$ws.Range("S7").Value = 0.4248704663212435
$ws.Range("B8").Value = 0.07434052757793765
$ws.Range("D8").Value = 0.01678657074340528
$ws.Range("E8").Value = 0.002398081534772182
$ws.Range("F8").Value = 0.03836930455635491
$ws.Range("J8").Value = 0.1199040767386091
$ws.Range("O8").Value = 0.01199040767386091
$ws.Range("Q8").Value = 0.1750599520383693
$ws.Range("R8").Value = 0.1294964028776978
$ws.Range("S8").Value = 0.4316546762589928
$ws.Range("B9").Value = 0.1089108910891089
$ws.Range("D9").Value = 0.0198019801980198
$ws.Range("F9").Value = 0.06435643564356436
$ws.Range("J9").Value = 0.103960396039604
$ws.Range("O9").Value = 0.0198019801980198
$ws.Range("Q9").Value = 0.1584158415841584
$ws.Range("R9").Value = 0.09405940594059406
$ws.Range("S9").Value = 0.4306930693069307
$ws.Range("B10").Value = 0.1086786551993745
$ws.Range("D10").Value = 0.01641907740422205
$ws.Range("F10").Value = 0.0727130570758405
$ws.Range("J10").Value = 0.1180609851446442
$ws.Range("O10").Value = 0.0218921032056294
$ws.Range("Q10").Value = 0.200938232994527
$ws.Range("R10").Value = 0.08522283033620015
$ws.Range("S10").Value = 0.3760750586395621
$ws.Range("G11").Value = 0.1389830508474576
$ws.Range("J11").Value = 0.07457627118644068
$ws.Range("K11").Value = 0.2067796610169491
$ws.Range("L11").Value = 0.559322033898305
$ws.Range("S11").Value = 0.02033898305084746
$ws.Range("G12").Value = 0.76
$ws.Range("J12").Value = 0.1714285714285714
$ws.Range("K12").Value = 0.01142857142857143
$ws.Range("L12").Value = 0.03428571428571429
$ws.Range("S12").Value = 0.02285714285714286
$ws.Range("G13").Value = 0.6944444444444444
$ws.Range("J13").Value = 0.2222222222222222
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("F15").Value = 0.04
$ws.Range("H15").Value = 0.1422222222222222
$ws.Range("I15").Value = 0.07111111111111111
$ws.Range("J15").Value = 0.3733333333333334
$ws.Range("K15").Value = 0.05333333333333334
$ws.Range("M15").Value = 0.004444444444444444
$ws.Range("N15").Value = 0.004444444444444444
$ws.Range("O15").Value = 0.06222222222222222
$ws.Range("S15").Value = 0.2488888888888889
$ws.Range("F16").Value = 0.02272727272727273
$ws.Range("H16").Value = 0.1875
$ws.Range("I16").Value = 0.08522727272727272
$ws.Range("J16").Value = 0.3863636363636364
$ws.Range("K16").Value = 0.1136363636363636
$ws.Range("M16").Value = 0.01136363636363636
$ws.Range("O16").Value = 0.03409090909090909
$ws.Range("S16").Value = 0.1590909090909091
$ws.Range("F17").Value = 0.02386634844868735
$ws.Range("H17").Value = 0.1742243436754177
$ws.Range("I17").Value = 0.1145584725536993
$ws.Range("J17").Value = 0.4033412887828162
$ws.Range("K17").Value = 0.09069212410501193
$ws.Range("M17").Value = 0.01670644391408115
$ws.Range("O17").Value = 0.05489260143198091
$ws.Range("S17").Value = 0.1217183770883055
$ws.Range("F18").Value = 0.04072398190045249
$ws.Range("H18").Value = 0.16289592760181
$ws.Range("I18").Value = 0.09954751131221719
$ws.Range("J18").Value = 0.3574660633484163
$ws.Range("K18").Value = 0.09049773755656108
$ws.Range("M18").Value = 0.02262443438914027
$ws.Range("O18").Value = 0.08144796380090498
$ws.Range("S18").Value = 0.1447963800904978
$ws.Range("F19").Value = 0.01566170712607674
$ws.Range("H19").Value = 0.1957713390759593
$ws.Range("I19").Value = 0.07909162098668755
$ws.Range("J19").Value = 0.3876272513703994
$ws.Range("K19").Value = 0.110415035238841
$ws.Range("M19").Value = 0.01957713390759593
$ws.Range("O19").Value = 0.07361002349256068
$ws.Range("S19").Value = 0.1182458888018794
